$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test1")

# Loadcell characterization update: replace formula-derived inclination
# angles with directly measured values, and refresh the initial voltage
# readings for both sensors.
$ws.Range("B16").Value = 82
$ws.Range("D16").Value = 1.81
$ws.Range("G16").Value = 97
$ws.Range("I16").Value = 3.43

# Drop the stray empty cell left over in the old layout.
$ws.Range("E18").ClearContents()

# Update the active selection to reflect the reviewed cell range.
$ws.Activate()
$ws.Range("C17:D17").Select()
